$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new data run row for sg_rr_100_027
$ws.Range("A33").Value = "sg_rr_100_027 2023-12-08 17-44-55.csv"
$ws.Range("B33").Value = 0.01
$ws.Range("C33").Value = 1000
$ws.Range("D33").Value = 5001
$ws.Range("E33").Value = 1530
$ws.Range("F33").Value = 1570
$ws.Range("G33").Value = 0.5
$ws.Range("H33").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I33").Value = 1.7
$ws.Range("J33").Value = 0.98153846153846003
$ws.Range("K33").Value = 0.0043858818636388196
$ws.Range("M33").Value = "prominence kept same, thought about increasing a little when trying to look roughly by eye for roughly biggest height span of noisy bit containing no peaks but decide to keep it same and adjust slightly if need be, as looks roughly right anyway, and wouldn't want to increase it unnecessarily if not needed."
$ws.Range("L33").Value = "yes (no double-counting but possible loss of last peak as it is on edge)"

# Update selection to reflect new active cell after entry (A34)
$ws.Range("A34").Select()
